$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/changed long text blocks
$objetivoGeral = "Objetivo GeralPermitir aos estudantes que compreendam os mecanismos de obten" + [char]0x00E7 + [char]0x00E3 + "o da influencia de diversos fatores (vari" + [char]0x00E1 + "veis independentes de um processo) sobre as vari" + [char]0x00E1 + "veis resposta (dependentes), atrav" + [char]0x00E9 + "s da an" + [char]0x00E1 + "lise multivariada.Objetivos Espec" + [char]0x00ED + "ficosSaber planejar e executar um experimento fatorial completo e fracionadoSaber analisar os resultados propondo a condi" + [char]0x00E7 + [char]0x00E3 + "o de melhor ajuste que otimiza os valores da vari" + [char]0x00E1 + "vel resposta na regi" + [char]0x00E3 + "o experimental estudadaDominar, pelo menos, um software comercial sobre o assuntoSaber modelar um processo, com base em dados emp" + [char]0x00ED + "ricos"

$docenteResp = "5840535 - Messias Borges Silva"

$programaResumido = "Introdu" + [char]0x00E7 + [char]0x00E3 + "o Experimenta" + [char]0x00E7 + [char]0x00E3 + "o convencional Experimentos Fatoriais completos Experimentos Fatoriais fracionados An" + [char]0x00E1 + "lise de vari" + [char]0x00E2 + "ncia Metodologia de superf" + [char]0x00ED + "cie de resposta M" + [char]0x00E9 + "todo de Taguchi"

$programa = "Introdu" + [char]0x00E7 + [char]0x00E3 + "o Experimenta" + [char]0x00E7 + [char]0x00E3 + "o convencional Experimentos Fatoriais completos 2k , Experimentos Fatoriais fracionados 2k-p, M" + [char]0x00E9 + "todo de Plackett Burman,  An" + [char]0x00E1 + "lise de vari" + [char]0x00E2 + "ncia Metodologia de superf" + [char]0x00ED + "cie de resposta, M" + [char]0x00E9 + "todo de Taguchi ."

$metodo = "2 provas escritas"

$criterio = "Ser" + [char]0x00E3 + "o avaliados os conte" + [char]0x00FA + "dos discutidos em sala e constantes da ementa do curso. MF = (0,40*P1 + 0,40*P2 + 0,20*TRAB), onde P1 e P2 s" + [char]0x00E3 + "o provas e TRAB " + [char]0x00E9 + " a nota m" + [char]0x00E9 + "dia de trabalhos e semin" + [char]0x00E1 + "rios."

$normaRecup = "Uma provas escrita com conte" + [char]0x00FA + "do de todo o semestre. NF = (MF + PR)/2, onde PR " + [char]0x00E9 + " uma prova de recupera" + [char]0x00E7 + [char]0x00E3 + "o"

$bibliografia = "1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 19912. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley & Sons, New York, 1978. 3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. 4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. 5. COX, D.R., Planning of Experiments, Wiley 1976. 6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. 7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013"

# 1) Row 10 (Objetivos:) B/C content replaced with the Portuguese objectives text
$ws.Range("B10:C10").Value = $objetivoGeral

# 2) Insert a new blank row at 13 (pushes existing rows 13.. down by one)
$ws.Rows.Item(13).Insert()

# 3) New row 13 gets the "Docentes responsaveis" value (moved content) in B/C, no A value.
#    Clear A13 completely (Insert() copies row-12's bold style into col A) and copy the
#    plain B/C styles down from row 9 so B13/C13 match the rest of the B/C columns.
$ws.Range("A13").Clear()
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13:C13").Value = $docenteResp

# 4) Row 14 (was 13, Programa resumido:) B/C: "Semestral" -> new syllabus summary text
$ws.Range("B14:C14").Value = $programaResumido

# 5) Row 16 (was 15, Programa:) B/C: "01/01/2018" -> new full syllabus text
$ws.Range("B16:C16").Value = $programa

# 6) Row 19 (was 18, Metodo:) B/C: "5840535 - Messias Borges Silva" -> "2 provas escritas"
$ws.Range("B19:C19").Value = $metodo

# 7) Row 20 (was 19, Criterio:) B/C: "2 provas escritas" -> evaluation criteria text
$ws.Range("B20:C20").Value = $criterio

# 8) Row 21 (was 20, Norma de recuperacao:) B/C: old criteria text -> makeup exam text
$ws.Range("B21:C21").Value = $normaRecup

# 9) Row 22 (was 21, Bibliografia:) B/C: old makeup text -> bibliography text
$ws.Range("B22:C22").Value = $bibliografia

Write-Output "edit complete"
